$wb = $excel.ActiveWorkbook

# --- Switch the active/selected sheet from "NewLoanInput" to "Repayment schedule" ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate()

# --- Insert a new (blank) column before column N, shifting the old N:P to O:Q ---
# Grab the column-width (in Excel's character units) of the column to the left (M)
# before inserting, so the new column can inherit that same width (Excel's default
# "insert" behaviour copies the format of the column to the left).
$leftWidth = $wsRepay.Columns("M").ColumnWidth
$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = $leftWidth

# --- Update the selection on the Repayment schedule sheet ---
$wsRepay.Range("R7").Select() | Out-Null
